# apex reworked for 4 TCDS links
# - Add a new "mem_base" parameter row (row 22) to both root-config sheets
#   (apex_ku15p_gty_root_config_2 and apex_ku15p_gth_root_config), mirroring
#   the existing "device_count" row's layout/formatting.
# - Make "apex_ku15p_gth_root_config" the active sheet/tab.
# - Update the selections left behind on each sheet to reflect the edit
#   (new row selected on the two config sheets; Instructions sheet keeps
#   its previous active cell after the range was referenced there too).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # apex_ku15p_gty_root_config_2
$ws2 = $wb.Worksheets.Item(2)   # apex_ku15p_gth_root_config
$ws3 = $wb.Worksheets.Item(3)   # Instructions

# --- Sheet 1 (apex_ku15p_gty_root_config_2): append new "mem_base" row ---
[void]$ws1.Range("A21:C21").Copy()
[void]$ws1.Range("A22:C22").PasteSpecial()
$ws1.Range("A22").Value = "mem_base"
$ws1.Range("B22").Value = "0x0"
$ws1.Range("C22").Value = "memory base for PCIe systems"
[void]$ws1.Range("A22:C22").Select()

# --- Sheet 2 (apex_ku15p_gth_root_config): append matching new row ---
[void]$ws2.Range("A21:C21").Copy()
[void]$ws2.Range("A22:C22").PasteSpecial()
$ws2.Range("A22").Value = "mem_base"
$ws2.Range("B22").Value = "0x0"
$ws2.Range("C22").Value = "memory base for PCIe systems"

# --- Instructions sheet: selection touches the copied range and the prior
#     active cell (C25) ---
[void]$ws3.Activate()
[void]$ws3.Range("C25,A22:C22").Select()

# --- Final state: apex_ku15p_gth_root_config is the active/selected tab ---
[void]$ws2.Activate()
[void]$ws2.Range("A22:C22").Select()
